# Slide 3 ("Roadmap") - Content Placeholder 2:
#   - "NOSQL " + "overview " -> single run "NOSQL overview "
#   - "fundamentals" -> "fundamentals (and tools)"
#   - drop the standalone "Gremlin" bullet entirely

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sp = $s.Shapes.Item(2)
$tr = $sp.TextFrame.TextRange

# Remove the "Gremlin" paragraph (including its trailing paragraph mark) first,
# from the tail of the text so earlier character offsets stay valid.
$gremlinPara = $tr.Characters(98, 8)
$gremlinPara.Delete()

# "fundamentals" -> "fundamentals (and tools)" (2nd run of 2nd paragraph)
$fundamentalsRun = $tr.Characters(23, 12)
$fundamentalsRun.Text = "fundamentals (and tools)"

# "NOSQL " + "overview " -> single run "NOSQL overview " (1st paragraph)
$titleRun = $tr.Characters(1, 15)
$titleRun.Text = "NOSQL overview "
